# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-12-01 07:18:53
#
# This script re-applies the "Recorded By" attendee-list re-ordering and the
# status/metric updates produced by the latest sync of the attendance report
# (the underlying email lists were re-ordered by the source system and one
# session flipped from "Pending" to "Not Recorded").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Recorded By" (column G) email list re-ordering ---------------------

$ws.Range("G2").Value = "Veronia.rafat@med.asu.edu.eg, System, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, gehanadel@med.asu.edu.eg"

$ws.Range("G3").Value = "asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, System, majorelle.magdy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"

$ws.Range("G4").Value = "asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, gehanadel@med.asu.edu.eg"

$ws.Range("G6").Value = "majorelle.magdy@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm"

$ws.Range("G7").Value = "Amera.a.saad@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg"

$ws.Range("G12").Value = "amira.m.ibrahim@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg"

$ws.Range("G15").Value = "mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"

$ws.Range("G27").Value = "hana.amr@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"

$ws.Range("G28").Value = "Aya_hamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"

# --- Status change: HISTOLOGY session 3 (row 11) is now "Not Recorded" ---

$ws.Range("I11").Value = "Not Recorded"

# --- Resulting metric updates (Class Statistics block) -------------------

$ws.Range("L7").Value = 2   # Missing Sessions: 1 -> 2
$ws.Range("L8").Value = 9   # Pending Sessions: 10 -> 9

# --- Resulting metric updates (Group Statistics block, row 15) -----------

$ws.Range("P15").Value = 2  # Missing: 1 -> 2
$ws.Range("Q15").Value = 9  # Pending: 10 -> 9

Write-Host "Applied attendance_reports sync updates"
